$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2163009404388715
$ws.Range("C2").Value = 0.5266457680250783
$ws.Range("J2").Value = 0.02507836990595611
$ws.Range("P2").Value = 0.1661442006269593
$ws.Range("S2").Value = 0.06583072100313479
$ws.Range("C3").Value = 0.02162162162162162
$ws.Range("J3").Value = 0.07567567567567568
$ws.Range("P3").Value = 0.7135135135135136
$ws.Range("S3").Value = 0.1891891891891892
$ws.Range("J4").Value = 0.1219512195121951
$ws.Range("P4").Value = 0.5365853658536586
$ws.Range("S4").Value = 0.3414634146341464
$ws.Range("B6").Value = 0.05494505494505494
$ws.Range("D6").Value = 0.01098901098901099
$ws.Range("F6").Value = 0.08241758241758242
$ws.Range("J6").Value = 0.3571428571428572
$ws.Range("O6").Value = 0.04395604395604396
$ws.Range("Q6").Value = 0.07142857142857142
$ws.Range("R6").Value = 0.08241758241758242
$ws.Range("S6").Value = 0.2967032967032967
$ws.Range("B7").Value = 0.09444444444444444
$ws.Range("D7").Value = 0.02777777777777778
$ws.Range("F7").Value = 0.03333333333333333
$ws.Range("J7").Value = 0.1666666666666667
$ws.Range("O7").Value = 0.01111111111111111
$ws.Range("Q7").Value = 0.1555555555555556
$ws.Range("R7").Value = 0.07777777777777778
$ws.Range("S7").Value = 0.4333333333333333
$ws.Range("B8").Value = 0.1006864988558352
$ws.Range("D8").Value = 0.02517162471395881
$ws.Range("F8").Value = 0.03432494279176201
$ws.Range("J8").Value = 0.1167048054919908
$ws.Range("O8").Value = 0.02517162471395881
$ws.Range("Q8").Value = 0.1807780320366133
$ws.Range("R8").Value = 0.07551487414187644
$ws.Range("S8").Value = 0.4416475972540046
$ws.Range("B9").Value = 0.125
$ws.Range("D9").Value = 0.0131578947368421
$ws.Range("F9").Value = 0.03947368421052631
$ws.Range("J9").Value = 0.1052631578947368
$ws.Range("O9").Value = 0.03289473684210526
$ws.Range("Q9").Value = 0.1578947368421053
$ws.Range("R9").Value = 0.1118421052631579
$ws.Range("S9").Value = 0.4144736842105263
$ws.Range("B10").Value = 0.1180124223602484
$ws.Range("D10").Value = 0.01552795031055901
$ws.Range("F10").Value = 0.05124223602484472
$ws.Range("J10").Value = 0.15527950310559
$ws.Range("O10").Value = 0.02018633540372671
$ws.Range("Q10").Value = 0.2204968944099379
$ws.Range("R10").Value = 0.06521739130434782
$ws.Range("S10").Value = 0.3540372670807453
$ws.Range("G11").Value = 0.1480263157894737
$ws.Range("J11").Value = 0.09539473684210527
$ws.Range("K11").Value = 0.2269736842105263
$ws.Range("L11").Value = 0.5098684210526315
$ws.Range("S11").Value = 0.01973684210526316
$ws.Range("G12").Value = 0.7278481012658228
$ws.Range("J12").Value = 0.2088607594936709
$ws.Range("K12").Value = 0.006329113924050633
$ws.Range("L12").Value = 0.03164556962025317
$ws.Range("S12").Value = 0.02531645569620253
$ws.Range("G13").Value = 0.6216216216216216
$ws.Range("J13").Value = 0.3513513513513514
$ws.Range("S13").Value = 0.02702702702702703
$ws.Range("F15").Value = 0.0205761316872428
$ws.Range("H15").Value = 0.1234567901234568
$ws.Range("I15").Value = 0.06172839506172839
$ws.Range("J15").Value = 0.3703703703703703
$ws.Range("K15").Value = 0.06584362139917696
$ws.Range("M15").Value = 0.00823045267489712
$ws.Range("O15").Value = 0.07407407407407407
$ws.Range("S15").Value = 0.2757201646090535
$ws.Range("F16").Value = 0.0101010101010101
$ws.Range("H16").Value = 0.1818181818181818
$ws.Range("I16").Value = 0.04545454545454546
$ws.Range("J16").Value = 0.4343434343434344
$ws.Range("K16").Value = 0.1414141414141414
$ws.Range("M16").Value = 0.005050505050505051
$ws.Range("O16").Value = 0.07575757575757576
$ws.Range("S16").Value = 0.1060606060606061
$ws.Range("F17").Value = 0.02293577981651376
$ws.Range("H17").Value = 0.1743119266055046
$ws.Range("I17").Value = 0.09174311926605505
$ws.Range("J17").Value = 0.4128440366972477
$ws.Range("K17").Value = 0.0871559633027523
$ws.Range("M17").Value = 0.02064220183486239
$ws.Range("O17").Value = 0.06192660550458716
$ws.Range("S17").Value = 0.1284403669724771
$ws.Range("F18").Value = 0.03636363636363636
$ws.Range("H18").Value = 0.2303030303030303
$ws.Range("I18").Value = 0.05454545454545454
$ws.Range("J18").Value = 0.3333333333333333
$ws.Range("K18").Value = 0.1090909090909091
$ws.Range("O18").Value = 0.09696969696969697
$ws.Range("S18").Value = 0.1393939393939394
$ws.Range("F19").Value = 0.02188006482982172
$ws.Range("H19").Value = 0.2090761750405186
$ws.Range("I19").Value = 0.06482982171799027
$ws.Range("J19").Value = 0.3638573743922204
$ws.Range("K19").Value = 0.1085899513776337
$ws.Range("M19").Value = 0.02106969205834684
$ws.Range("N19").Value = 0.0008103727714748784
$ws.Range("O19").Value = 0.07617504051863858
$ws.Range("S19").Value = 0.1337115072933549
